$wb = $excel.ActiveWorkbook

# --- Rename existing sheet, add the new summary sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Stock log"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Portfolio Summary"

# Match the sheet-level properties used throughout this workbook (outline
# direction + page margins), mirroring the "Stock log" sheet.
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# --- Sheet "Stock log": insert two new columns (Cost Basis, Market Value) ---
# before the existing "Capital Gains" column, shifting H,I,J -> J,K,L
$ws1.Range("H1:I1").EntireColumn.Insert()

# Header row
$ws1.Range("H1").Value = "Cost Basis"
$ws1.Range("I1").Value = "Market Value"
$ws1.Range("J1").Value = "Capital Gains"
$ws1.Range("K1").Value = "Dividends Paid"
$ws1.Range("L1").Value = "Total Return"

# Updated Market Price values
$ws1.Range("E2").Value = 21.53
$ws1.Range("E3").Value = 21.53

# Row 2 (AY, 687 shares)
$ws1.Range("H2").Value = 16178.85
$ws1.Range("I2").Value = 14791.11
$ws1.Range("J2").Value = -8.58
$ws1.Range("K2").Value = 611.4300000000001
$ws1.Range("L2").Value = -4.99

# Row 3 (AY, 400 shares)
$ws1.Range("H3").Value = 9864
$ws1.Range("I3").Value = 8612
$ws1.Range("J3").Value = -12.69
$ws1.Range("K3").Value = 356
$ws1.Range("L3").Value = -9.42

# Row 4 (SCHD, 983 shares)
$ws1.Range("H4").Value = 74363.95000000001
$ws1.Range("I4").Value = 74531.06
$ws1.Range("J4").Value = 0.22
$ws1.Range("K4").Value = 1376.2
$ws1.Range("L4").Value = 2.11

# --- Sheet "Portfolio Summary": pivot-style summary, one row per ticker ---
$ws2.Range("A1").Value = "Number of Shares"
$ws2.Range("B1").Value = "Cost Basis"
$ws2.Range("C1").Value = "Market Value"
$ws2.Range("D1").Value = "Dividends Paid"
$ws2.Range("E1").Value = "Average price paid, USD"
$ws2.Range("F1").Value = "Capital Gains, %"
$ws2.Range("G1").Value = "Total Return, %"

# Match the header style used on the "Stock log" sheet (bold, bordered,
# centered/top-aligned) by copying its formatting over.
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122) | Out-Null

# AY (Atlantica Sustainable Infrastructure) - combines rows 2 and 3 from Stock log
$ws2.Range("A2").Value = 1087
$ws2.Range("B2").Value = 26042.85
$ws2.Range("C2").Value = 23403.11
$ws2.Range("D2").Value = 967.4300000000001
$ws2.Range("E2").Value = 23.96
$ws2.Range("F2").Value = -10.14
$ws2.Range("G2").Value = -6.68

# SCHD (Schwab US Dividend Equity ETF)
$ws2.Range("A3").Value = 983
$ws2.Range("B3").Value = 74363.95000000001
$ws2.Range("C3").Value = 74531.06
$ws2.Range("D3").Value = 1376.2
$ws2.Range("E3").Value = 75.65000000000001
$ws2.Range("F3").Value = 0.22
$ws2.Range("G3").Value = 2.11

# Keep "Stock log" as the active/selected sheet (matches original activeTab="0")
$ws1.Activate()
